$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.905.13"
$ws.Range("E2").Value = "  -0.07%  "
$ws.Range("D3").Value = "1.814.73"
$ws.Range("E3").Value = "  +1.50%  "
$ws.Range("D4").Value = "'1.005"
$ws.Range("E4").Value = "  -0.75%  "
$ws.Range("D5").Value = "'311.23"
$ws.Range("E5").Value = "  -0.05%  "
$ws.Range("D6").Value = "'1.004"
$ws.Range("E6").Value = "  -0.37%  "
$ws.Range("D7").Value = "'0.4289"
$ws.Range("E7").Value = "  +1.41%  "
$ws.Range("D8").Value = "'0.3686"
$ws.Range("E8").Value = "  +2.40%  "
$ws.Range("D9").Value = "'0.07254"
$ws.Range("E9").Value = "  +1.19%  "
$ws.Range("D10").Value = "'0.8616"
$ws.Range("E10").Value = "  +2.54%  "
$ws.Range("E11").Value = "  +4.01%  "
$ws.Range("D12").Value = "2.010.51"
$ws.Range("E12").Value = "  +8.81%  "
$ws.Range("D13").Value = "'6.640"
$ws.Range("E13").Value = "  +4.55%  "
$ws.Range("D14").Value = "'5.401"
$ws.Range("E14").Value = "  +2.81%  "
$ws.Range("D15").Value = "'0.06901"
$ws.Range("E15").Value = "  +1.02%  "
$ws.Range("D16").Value = "'80.66"
$ws.Range("E16").Value = "  +1.10%  "
$ws.Range("E17").Value = "  -0.78%  "
$ws.Range("D18").Value = "'0.000008928"
$ws.Range("E18").Value = "  +2.64%  "
$ws.Range("E19").Value = "  -0.43%  "
$ws.Range("E20").Value = "  +1.83%  "
$ws.Range("D21").Value = "26.958.20"
$ws.Range("E21").Value = "  -0.63%  "
$ws.Range("D22").Value = "'5.182"
$ws.Range("E22").Value = "  +2.74%  "
$ws.Range("D23").Value = "'11.07"
$ws.Range("E23").Value = "  +0.17%  "
$ws.Range("D24").Value = "2.228.76"
$ws.Range("E24").Value = "  +8.04%  "
$ws.Range("D25").Value = "'153.78"
$ws.Range("E25").Value = "  +0.54%  "
$ws.Range("D26").Value = "'1.885"
$ws.Range("E26").Value = "  -3.33%  "
$ws.Range("D27").Value = "'18.24"
$ws.Range("E27").Value = "  +0.36%  "
$ws.Range("D28").Value = "'5.196"
$ws.Range("E28").Value = "  +3.78%  "
$ws.Range("D29").Value = "'115.02"
$ws.Range("E29").Value = "  +0.22%  "
$ws.Range("D30").Value = "'1.874"
$ws.Range("E30").Value = "  +16.13%  "
$ws.Range("D31").Value = "'0.08950"
$ws.Range("E31").Value = "  +0.23%  "
$ws.Range("D32").Value = "'0.7424"
$ws.Range("E32").Value = "  +2.82%  "
$ws.Range("E33").Value = "  +7.59%  "
$ws.Range("D34").Value = "'4.419"
$ws.Range("E34").Value = "  +2.40%  "
$ws.Range("D35").Value = "'2.801"
$ws.Range("E35").Value = "  -1.67%  "
$ws.Range("D36").Value = "'1.009"
$ws.Range("E36").Value = "  -0.08%  "
$ws.Range("E37").Value = "  +4.08%  "
$ws.Range("D38").Value = "'0.05211"
$ws.Range("E38").Value = "  +2.51%  "
$ws.Range("E39").Value = "  +1.80%  "
$ws.Range("D40").Value = "'0.5067"
$ws.Range("E40").Value = "  +2.69%  "
$ws.Range("B41").Value = "Algorand"
$ws.Range("C41").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D41").Value = "'0.1643"
$ws.Range("E41").Value = "  +1.89%  "
$ws.Range("B42").Value = "MXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D42").Value = "'2.726"
$ws.Range("E42").Value = "  +9.33%  "
$ws.Range("D43").Value = "'6.421"
$ws.Range("E43").Value = "  +7.66%  "
$ws.Range("D44").Value = "'8.234"
$ws.Range("E44").Value = "  +3.90%  "
$ws.Range("D45").Value = "'106.68"
$ws.Range("E45").Value = "  +2.20%  "
$ws.Range("D46").Value = "'10.42"
$ws.Range("E46").Value = "  +3.65%  "
$ws.Range("D47").Value = "'1.005"
$ws.Range("E47").Value = "  -0.32%  "
$ws.Range("D48").Value = "'1.653"
$ws.Range("E48").Value = "  +4.93%  "
$ws.Range("D49").Value = "'0.06294"
$ws.Range("E49").Value = "  +0.15%  "
$ws.Range("D50").Value = "'0.4551"
$ws.Range("E50").Value = "  +1.67%  "
$ws.Range("D51").Value = "'1.797"
$ws.Range("E51").Value = "  +5.62%  "
